# The author deleted the post row that used to sit at row 814
# ("「サイのざらざらした皮を利用して自分の体を掻くネコ」" / rhino-skin cat post).
# Deleting the entire row shifts every row below it up by one and shrinks
# the sheet's used range from A1:C882 to A1:C881, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(814).Delete()
